# Pedido 69134d11b9c1d30b15fabdc3
# -----------------------------------------------------------------------
# The sheet "Productos" currently ends at row 4 (A1:N4). This change adds
# a brand-new row 5 that duplicates row 4 verbatim (same "Test Ringover
# (NO TOCAR)" product line), growing the used range to A1:N5.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 4
$newRow = $lastDataRow + 1
$lastCol = 14   # column N ("Fecha de venta")

# Copy the whole last row and insert the copy as a new row right below
# it (mirrors Excel's "Copy row -> Insert Copied Cells" workflow). This
# keeps every cell's original data type (numbers stay numbers, text
# stays text) instead of re-typing values through .Value/.Value2.
$ws.Rows.Item($lastDataRow).Copy()
$ws.Rows.Item($newRow).Insert()

# A few columns on the source row (Optimizador, Unidades Optimizador,
# Cargador VE) are blank-but-present empty-text placeholders rather than
# truly empty cells. Those don't survive the row copy/insert as empty
# text cells, so detect and recreate them explicitly on the new row too.
for ($c = 1; $c -le $lastCol; $c++) {
    $srcCell = $ws.Cells.Item($lastDataRow, $c)
    $srcVal = $srcCell.Value2
    $isEmptyText = ($srcVal -ne $null) -and ($srcVal.GetType().Name -eq "String") -and ($srcVal -eq "")
    if ($isEmptyText) {
        $dstCell = $ws.Cells.Item($newRow, $c)
        $dstCell.Formula = "'"
        $dstCell.ClearFormats()
    }
}
